$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1490
$ws.Range("I40").Value = 1071.4286
$ws.Range("J40").Value = 2466.6667
$ws.Range("K40").Value = 1071.4286
$ws.Range("L40").Value = 2466.6667
$ws.Range("M40").Value = -896.4286
$ws.Range("N40").Value = -2816.6667

$ws.Range("H62").Value = 3497.738
$ws.Range("I62").Value = 1455.3125
$ws.Range("J62").Value = 4754.615
$ws.Range("K62").Value = 1455.3125
$ws.Range("L62").Value = 4754.615
$ws.Range("M62").Value = -831.3125
$ws.Range("N62").Value = -6002.615

$ws.Range("H65").Value = 3497.738
$ws.Range("I65").Value = 1455.3125
$ws.Range("J65").Value = 4754.615
$ws.Range("K65").Value = 7276.5625
$ws.Range("L65").Value = 23773.075
$ws.Range("M65").Value = -4156.5625
$ws.Range("N65").Value = -30013.075

$ws.Range("H125").Value = 2306.077
$ws.Range("I125").Value = 1097.5
$ws.Range("J125").Value = 4239.8
$ws.Range("K125").Value = 9877.5
$ws.Range("L125").Value = 38158.2
$ws.Range("M125").Value = -7417.5
$ws.Range("N125").Value = -43078.2

$ws.Range("H132").Value = 10564.091
$ws.Range("I132").Value = 12133.889
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 36401.667
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -33871.667
$ws.Range("N132").Value = -15560

$ws.Range("H135").Value = 12327139
$ws.Range("I135").Value = 708.125
$ws.Range("J135").Value = 27897366
$ws.Range("K135").Value = 6373.125
$ws.Range("L135").Value = 251076294
$ws.Range("M135").Value = -3838.125
$ws.Range("N135").Value = -251081364

$ws.Range("H138").Value = 1958.0541
$ws.Range("I138").Value = 1555.4814
$ws.Range("J138").Value = 3045
$ws.Range("K138").Value = 4666.4442
$ws.Range("L138").Value = 9135
$ws.Range("M138").Value = 473.5558000000001
$ws.Range("N138").Value = -19415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8386.286
$ws.Range("J37").Value = 11219.8
$ws.Range("L37").Value = 11219.8
$ws.Range("N37").Value = -11765.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13031.25
$ws.Range("I82").Value = 6929.3335
$ws.Range("J82").Value = 20876.572
$ws.Range("K82").Value = 6929.3335
$ws.Range("L82").Value = 20876.572
$ws.Range("M82").Value = -6546.3335
$ws.Range("N82").Value = -21642.572

$ws.Range("H85").Value = 13031.25
$ws.Range("I85").Value = 6929.3335
$ws.Range("J85").Value = 20876.572
$ws.Range("K85").Value = 6929.3335
$ws.Range("L85").Value = 20876.572
$ws.Range("M85").Value = -5603.3335
$ws.Range("N85").Value = -23528.572

$ws.Range("H134").Value = 6935448.5
$ws.Range("I134").Value = 8744391
$ws.Range("J134").Value = 1169
$ws.Range("K134").Value = 26233173
$ws.Range("L134").Value = 3507
$ws.Range("M134").Value = -26230638
$ws.Range("N134").Value = -8577

$ws.Range("H135").Value = 37490.59
$ws.Range("J135").Value = 37490.59
$ws.Range("L135").Value = 37490.59
$ws.Range("N135").Value = -47630.59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10218
$ws.Range("J50").Value = 11680.2
$ws.Range("L50").Value = 11680.2
$ws.Range("N50").Value = -12930.2

$ws.Range("H51").Value = 17870
$ws.Range("I51").Value = 12490
$ws.Range("J51").Value = 23250
$ws.Range("K51").Value = 12490
$ws.Range("L51").Value = 23250
$ws.Range("M51").Value = -11754
$ws.Range("N51").Value = -24722

$ws.Range("H61").Value = 17870
$ws.Range("I61").Value = 12490
$ws.Range("J61").Value = 23250
$ws.Range("K61").Value = 12490
$ws.Range("L61").Value = 23250
$ws.Range("M61").Value = -12142
$ws.Range("N61").Value = -23946

$ws.Range("H68").Value = 14000
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 20000
$ws.Range("M68").Value = -9251
$ws.Range("N68").Value = -21498

$ws.Range("H71").Value = 14000
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 60000
$ws.Range("M71").Value = -26256
$ws.Range("N71").Value = -67488

$ws.Range("H86").Value = 5084.278
$ws.Range("I86").Value = 6668.8335
$ws.Range("J86").Value = 4292
$ws.Range("K86").Value = 6668.8335
$ws.Range("L86").Value = 4292
$ws.Range("M86").Value = -5545.8335
$ws.Range("N86").Value = -6538

$ws.Range("H89").Value = 5084.278
$ws.Range("I89").Value = 6668.8335
$ws.Range("J89").Value = 4292
$ws.Range("K89").Value = 33344.1675
$ws.Range("L89").Value = 21460
$ws.Range("M89").Value = -27728.1675
$ws.Range("N89").Value = -32692

$ws.Range("H132").Value = 4520.9165
$ws.Range("I132").Value = 5904
$ws.Range("J132").Value = 3137.8333
$ws.Range("K132").Value = 17712
$ws.Range("L132").Value = 9413.499899999999
$ws.Range("M132").Value = -15182
$ws.Range("N132").Value = -14473.4999

$ws.Range("H134").Value = 3130.2
$ws.Range("I134").Value = 3130.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9390.599999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6855.599999999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5701
$ws.Range("J101").Value = 6596
$ws.Range("L101").Value = 19788
$ws.Range("N101").Value = -24656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 1571.375
$ws.Range("I40").Value = 1178.2
$ws.Range("J40").Value = 2226.6667
$ws.Range("K40").Value = 1178.2
$ws.Range("L40").Value = 2226.6667
$ws.Range("M40").Value = -1042.2
$ws.Range("N40").Value = -2498.6667

$ws.Range("H122").Value = 1904.1052
$ws.Range("I122").Value = 1927.7222
$ws.Range("J122").Value = 1479
$ws.Range("K122").Value = 5783.1666
$ws.Range("L122").Value = 4437
$ws.Range("M122").Value = -3333.1666
$ws.Range("N122").Value = -9337

